$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 462.7143
$ws.Range("I2").Value = 462.7143
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 462.7143
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -349.7143
$ws.Range("N2").ClearContents()
$ws.Range("H33").Value = 778.9231
$ws.Range("I33").Value = 831.3333
$ws.Range("K33").Value = 831.3333
$ws.Range("M33").Value = -602.3333
$ws.Range("H43").Value = 6327.6665
$ws.Range("I43").Value = 6327.6665
$ws.Range("K43").Value = 6327.6665
$ws.Range("M43").Value = -6258.6665
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 8101
$ws.Range("I26").Value = 8101
$ws.Range("K26").Value = 8101
$ws.Range("M26").Value = -7771
$ws.Range("H122").Value = 1204.8889
$ws.Range("I122").Value = 1171.6666
$ws.Range("J122").Value = 1271.3334
$ws.Range("K122").Value = 3514.9998
$ws.Range("L122").Value = 3814.0002
$ws.Range("M122").Value = -1064.9998
$ws.Range("N122").Value = -8714.0002
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -105060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 1000
$ws.Range("K23").Value = 1000
$ws.Range("M23").Value = -760
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -808
$ws.Range("H58").Value = 6401.1665
$ws.Range("I58").Value = 851.75
$ws.Range("K58").Value = 851.75
$ws.Range("M58").Value = -648.75
$ws.Range("H69").Value = 13035.4
$ws.Range("I69").Value = 8294.25
$ws.Range("K69").Value = 8294.25
$ws.Range("M69").Value = -7545.25
$ws.Range("H72").Value = 13035.4
$ws.Range("I72").Value = 8294.25
$ws.Range("K72").Value = 24882.75
$ws.Range("M72").Value = -21138.75
$ws.Range("H82").Value = 82498.75
$ws.Range("J82").Value = 82498.75
$ws.Range("L82").Value = 82498.75
$ws.Range("N82").Value = -83220.75
$ws.Range("H85").Value = 82498.75
$ws.Range("J85").Value = 82498.75
$ws.Range("L85").Value = 82498.75
$ws.Range("N85").Value = -84994.75
$ws.Range("H136").Value = 6401.1665
$ws.Range("I136").Value = 851.75
$ws.Range("K136").Value = 2555.25
$ws.Range("M136").Value = -5.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 39
$ws.Range("I7").Value = 22.625
$ws.Range("K7").Value = 67.875
$ws.Range("M7").Value = 44.125
$ws.Range("H12").Value = 36
$ws.Range("J12").Value = 28.833334
$ws.Range("L12").Value = 86.50000199999999
$ws.Range("N12").Value = -432.500002
$ws.Range("H23").Value = 218.2
$ws.Range("I23").Value = 247.75
$ws.Range("K23").Value = 743.25
$ws.Range("M23").Value = -508.25
$ws.Range("H68").Value = 1099.1428
$ws.Range("I68").Value = 1224
$ws.Range("K68").Value = 3672
$ws.Range("M68").Value = -2861
$ws.Range("H71").Value = 1099.1428
$ws.Range("I71").Value = 1224
$ws.Range("K71").Value = 11016
$ws.Range("M71").Value = -6960
$ws.Range("H117").Value = 4199.8
$ws.Range("I117").Value = 2999.6667
$ws.Range("J117").Value = 6000
$ws.Range("K117").Value = 8999.000100000001
$ws.Range("L117").Value = 18000
$ws.Range("M117").Value = -5557.000100000001
$ws.Range("N117").Value = -24884
$ws.Range("H121").Value = 1246.8572
$ws.Range("I121").Value = 450
$ws.Range("J121").Value = 1379.6666
$ws.Range("K121").Value = 1350
$ws.Range("L121").Value = 4138.9998
$ws.Range("M121").Value = -40
$ws.Range("N121").Value = -6758.9998
$ws.Range("H131").Value = 3660.6365
$ws.Range("I131").Value = 3963
$ws.Range("K131").Value = 11889
$ws.Range("M131").Value = -6849

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 25000
$ws.Range("J15").Value = 25000
$ws.Range("L15").Value = 25000
$ws.Range("N15").Value = -25576
$ws.Range("H62").Value = 7000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 7000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H81").Value = 25000
$ws.Range("J81").Value = 25000
$ws.Range("L81").Value = 25000
$ws.Range("N81").Value = -26996
$ws.Range("H84").Value = 25000
$ws.Range("J84").Value = 25000
$ws.Range("L84").Value = 75000
$ws.Range("N84").Value = -84984
$ws.Range("H113").Value = 12503
$ws.Range("I113").Value = 12503
$ws.Range("K113").Value = 12503
$ws.Range("M113").Value = -10333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1368.8889
$ws.Range("I22").Value = 902.8570999999999
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 902.8570999999999
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -607.8570999999999
$ws.Range("N22").Value = -3590
$ws.Range("H26").Value = 2000
$ws.Range("I26").Value = 2000
$ws.Range("K26").Value = 2000
$ws.Range("M26").Value = -1705
$ws.Range("H27").Value = 1368.8889
$ws.Range("I27").Value = 902.8570999999999
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 902.8570999999999
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -795.8570999999999
$ws.Range("N27").Value = -3214
$ws.Range("H46").Value = 6327.5
$ws.Range("I46").Value = 5798.3335
$ws.Range("K46").Value = 5798.3335
$ws.Range("M46").Value = -5610.3335
$ws.Range("H63").Value = 26692.334
$ws.Range("I63").Value = 20077
$ws.Range("J63").Value = 30000
$ws.Range("K63").Value = 20077
$ws.Range("L63").Value = 30000
$ws.Range("M63").Value = -19328
$ws.Range("N63").Value = -31498
$ws.Range("H66").Value = 26692.334
$ws.Range("I66").Value = 20077
$ws.Range("J66").Value = 30000
$ws.Range("K66").Value = 60231
$ws.Range("L66").Value = 90000
$ws.Range("M66").Value = -56487
$ws.Range("N66").Value = -97488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 100000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 100000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 100000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -101248
$ws.Range("H65").Value = 100000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 100000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 500000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -506240
$ws.Range("H122").Value = 1697.7858
$ws.Range("I122").Value = 1232.3334
$ws.Range("J122").Value = 2046.875
$ws.Range("K122").Value = 3697.0002
$ws.Range("L122").Value = 6140.625
$ws.Range("M122").Value = -1247.0002
$ws.Range("N122").Value = -11040.625
